$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 formulas: (Bn/Bn2)*100 instead of Bn/Bn2
$cols = @("B","C","D","E","F","G","H","I","J")
foreach ($col in $cols) {
    $cell = $ws.Range($col + "4")
    $cell.Formula = "=(" + $col + "3/" + $col + "2)*100"
}

# Delete row 5 entirely (Season/Summer/Winter/Spring data)
$ws.Rows("5").Delete()

# Set row 4 height
$ws.Rows("4").RowHeight = 16.2

# Update selection to I11
$ws.Range("I11").Select()
